# Weekly price update: insert a new observation row for Haba (Vega Monumental
# Concepción) right after the existing row 45, pushing the prior rows 46-58
# down to 47-59 (their data is preserved, only their row numbers change).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 46; existing rows 46..58 shift to 47..59.
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with this week's record.
$ws.Cells.Item(46, 1).Value = 11
$ws.Cells.Item(46, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(46, 3).Value = "Bíobío"
$ws.Cells.Item(46, 4).Value = 45120
$ws.Cells.Item(46, 5).Value = 8
$ws.Cells.Item(46, 6).Value = 100112026
$ws.Cells.Item(46, 7).Value = "Haba"
$ws.Cells.Item(46, 8).Value = "Sin especificar"
$ws.Cells.Item(46, 9).Value = "Primera"
$ws.Cells.Item(46, 10).Value = 100
$ws.Cells.Item(46, 11).Value = 16000
$ws.Cells.Item(46, 12).Value = 17000
$ws.Cells.Item(46, 13).Value = 16500
$ws.Cells.Item(46, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(46, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(46, 16).Value = 660
$ws.Cells.Item(46, 17).Value = 25
$ws.Cells.Item(46, 18).Value = "Hortaliza"
